# Applies updated transition-probability values to Sheet1 (UT Arlington_A matrix)
# per "changes to team matrices from games pulled march 7".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1758241758241758
$ws.Range("C2").Value = 0.6043956043956044
$ws.Range("J2").Value = 0.008241758241758242
$ws.Range("P2").Value = 0.1346153846153846
$ws.Range("S2").Value = 0.07692307692307693

# Row 3
$ws.Range("B3").Value = 0.02575107296137339
$ws.Range("C3").Value = 0.02145922746781116
$ws.Range("J3").Value = 0.02145922746781116
$ws.Range("P3").Value = 0.7510729613733905
$ws.Range("S3").Value = 0.1802575107296137

# Row 4
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("P4").Value = 0.5813953488372093
$ws.Range("S4").Value = 0.3720930232558139

# Row 6
$ws.Range("B6").Value = 0.05714285714285714
$ws.Range("D6").Value = 0.01071428571428571
$ws.Range("F6").Value = 0.06071428571428571
$ws.Range("J6").Value = 0.2714285714285714
$ws.Range("O6").Value = 0.05
$ws.Range("Q6").Value = 0.1535714285714286
$ws.Range("R6").Value = 0.05714285714285714
$ws.Range("S6").Value = 0.3392857142857143

# Row 7
$ws.Range("B7").Value = 0.1450980392156863
$ws.Range("D7").Value = 0.01568627450980392
$ws.Range("F7").Value = 0.06274509803921569
$ws.Range("J7").Value = 0.1333333333333333
$ws.Range("O7").Value = 0.007843137254901961
$ws.Range("Q7").Value = 0.1254901960784314
$ws.Range("R7").Value = 0.07450980392156863
$ws.Range("S7").Value = 0.4352941176470588

# Row 8
$ws.Range("B8").Value = 0.1191489361702128
$ws.Range("D8").Value = 0.01702127659574468
$ws.Range("F8").Value = 0.07872340425531915
$ws.Range("J8").Value = 0.1425531914893617
$ws.Range("O8").Value = 0.02978723404255319
$ws.Range("Q8").Value = 0.1340425531914894
$ws.Range("R8").Value = 0.09361702127659574
$ws.Range("S8").Value = 0.3851063829787234

# Row 9
$ws.Range("B9").Value = 0.08018867924528301
$ws.Range("D9").Value = 0.01886792452830189
$ws.Range("F9").Value = 0.1132075471698113
$ws.Range("J9").Value = 0.1320754716981132
$ws.Range("O9").Value = 0.02830188679245283
$ws.Range("Q9").Value = 0.1367924528301887
$ws.Range("R9").Value = 0.08018867924528301
$ws.Range("S9").Value = 0.4103773584905661

# Row 10
$ws.Range("B10").Value = 0.1172035640849897
$ws.Range("D10").Value = 0.0205620287868403
$ws.Range("E10").Value = 0.0006854009595613434
$ws.Range("F10").Value = 0.06716929403701165
$ws.Range("J10").Value = 0.1267991775188485
$ws.Range("O10").Value = 0.0205620287868403
$ws.Range("Q10").Value = 0.1788896504455106
$ws.Range("R10").Value = 0.09801233721727211
$ws.Range("S10").Value = 0.3701165181631254

# Row 11
$ws.Range("G11").Value = 0.1381074168797954
$ws.Range("J11").Value = 0.09974424552429667
$ws.Range("K11").Value = 0.1790281329923274
$ws.Range("L11").Value = 0.5677749360613811
$ws.Range("S11").Value = 0.01534526854219949

# Row 12
$ws.Range("G12").Value = 0.7603305785123967
$ws.Range("J12").Value = 0.1694214876033058
$ws.Range("K12").Value = 0.008264462809917356
$ws.Range("L12").Value = 0.02066115702479339
$ws.Range("S12").Value = 0.04132231404958678

# Row 15
$ws.Range("F15").Value = 0.03050847457627119
$ws.Range("H15").Value = 0.1152542372881356
$ws.Range("I15").Value = 0.08135593220338982
$ws.Range("J15").Value = 0.3694915254237288
$ws.Range("K15").Value = 0.0711864406779661
$ws.Range("M15").Value = 0.003389830508474576
$ws.Range("O15").Value = 0.07457627118644068
$ws.Range("S15").Value = 0.2542372881355932

# Row 16
$ws.Range("F16").Value = 0.03404255319148936
$ws.Range("H16").Value = 0.1787234042553192
$ws.Range("I16").Value = 0.04680851063829787
$ws.Range("J16").Value = 0.3829787234042553
$ws.Range("K16").Value = 0.1659574468085106
$ws.Range("M16").Value = 0.02553191489361702
$ws.Range("O16").Value = 0.02978723404255319
$ws.Range("S16").Value = 0.1361702127659574

# Row 17
$ws.Range("F17").Value = 0.0185614849187935
$ws.Range("H17").Value = 0.1647331786542924
$ws.Range("I17").Value = 0.08120649651972157
$ws.Range("J17").Value = 0.4083526682134571
$ws.Range("K17").Value = 0.1183294663573086
$ws.Range("M17").Value = 0.0185614849187935
$ws.Range("O17").Value = 0.06728538283062645
$ws.Range("S17").Value = 0.122969837587007

# Row 18
$ws.Range("F18").Value = 0.01687763713080169
$ws.Range("H18").Value = 0.1729957805907173
$ws.Range("I18").Value = 0.1012658227848101
$ws.Range("J18").Value = 0.4303797468354431
$ws.Range("K18").Value = 0.08438818565400844
$ws.Range("M18").Value = 0.02109704641350211
$ws.Range("N18").Value = 0.004219409282700422
$ws.Range("O18").Value = 0.05907172995780591
$ws.Range("S18").Value = 0.109704641350211

# Row 19
$ws.Range("F19").Value = 0.02517006802721088
$ws.Range("H19").Value = 0.1870748299319728
$ws.Range("I19").Value = 0.07278911564625851
$ws.Range("J19").Value = 0.3591836734693877
$ws.Range("K19").Value = 0.1231292517006803
$ws.Range("M19").Value = 0.01496598639455782
$ws.Range("N19").Value = 0.002040816326530612
$ws.Range("O19").Value = 0.08367346938775511
$ws.Range("S19").Value = 0.1319727891156462
